$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates apply identically to the "展览" sheet and
# the "全部类型" sheet (their data is mirrored).
$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    3  = 2711
    4  = 587
    5  = 90
    7  = 949
    8  = 10
    9  = 16
    11 = 64
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
